# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that the "System" token is moved to the front of the comma-separated list,
# while preserving the relative order of the remaining tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $splitResult = $val -split ", "

    # Copy the split result into a plain array because indexing directly
    # into the result of -split is unreliable in this runtime.
    $parts = New-Object 'object[]' $splitResult.Count
    $i = 0
    foreach ($p in $splitResult) {
        $parts[$i] = $p
        $i = $i + 1
    }

    if ($parts.Count -le 1) {
        continue
    }

    if (([string]$parts[0]).Equals("System")) {
        continue
    }

    $sysIndex = -1
    for ($j = 0; $j -lt $parts.Count; $j++) {
        if (([string]$parts[$j]).Equals("System")) {
            $sysIndex = $j
        }
    }

    if ($sysIndex -lt 0) {
        continue
    }

    $newParts = New-Object 'object[]' $parts.Count
    $newParts[0] = "System"
    $k = 1
    for ($j = 0; $j -lt $parts.Count; $j++) {
        if ($j -ne $sysIndex) {
            $newParts[$k] = $parts[$j]
            $k = $k + 1
        }
    }

    $newVal = [string]::Join(", ", $newParts)
    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
